$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st worksheet) - update F column "想去人数" values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1886
$ws1.Range("F3").Value = 499
$ws1.Range("F5").Value = 177
$ws1.Range("F6").Value = 2629
$ws1.Range("F7").Value = 174
$ws1.Range("F9").Value = 176
$ws1.Range("F10").Value = 1550
$ws1.Range("F11").Value = 538
$ws1.Range("F12").Value = 45
$ws1.Range("F14").Value = 233
$ws1.Range("F17").Value = 213
$ws1.Range("F21").Value = 189
$ws1.Range("F22").Value = 64
$ws1.Range("F23").Value = 1686
$ws1.Range("F24").Value = 36
$ws1.Range("F25").Value = 413
$ws1.Range("F26").Value = 20
$ws1.Range("F27").Value = 570
$ws1.Range("F28").Value = 211
$ws1.Range("F30").Value = 429

# Sheet "全部类型" (4th worksheet) - update F column "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1886
$ws4.Range("F4").Value = 499
$ws4.Range("F6").Value = 177
$ws4.Range("F7").Value = 2629
$ws4.Range("F8").Value = 174
$ws4.Range("F10").Value = 176
$ws4.Range("F11").Value = 1550
$ws4.Range("F12").Value = 538
$ws4.Range("F13").Value = 45
$ws4.Range("F15").Value = 233
$ws4.Range("F18").Value = 213
$ws4.Range("F22").Value = 189
$ws4.Range("F23").Value = 64
$ws4.Range("F24").Value = 1686
$ws4.Range("F25").Value = 36
$ws4.Range("F26").Value = 413
$ws4.Range("F27").Value = 20
$ws4.Range("F28").Value = 570
$ws4.Range("F29").Value = 211
$ws4.Range("F31").Value = 429
